# cosine_clustered_item_description.xlsx — "sent a copy to Nidhi"
#
# 1) The yellow highlight (and the blank Port_of_Shipment "G" cell that
#    goes with it) moves from row 3 to row 44.
# 2) A handful of rows have their USD-converted columns (T/U/V) refreshed
#    to slightly different values (an updated FX rate was applied).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Move the highlight from row 3 to row 44 -----------------------
# Copy row 3's formatting (fill, and the empty G cell that belongs to the
# highlighted block) onto row 44, then strip row 3 back down to plain/
# un-highlighted, fully removing its now-unused G3 cell.
$ws.Range("A3:V3").Copy()
$ws.Range("A44:V44").PasteSpecial(-4122)
$ws.Range("A3:V3").ClearFormats()
$ws.Range("G3").Clear()

# --- 2) Refresh the USD conversion columns -----------------------------
$ws.Range("T12").Value = 2.2005
$ws.Range("U12").Value = 2860.6286
$ws.Range("V12").Value = 3.4818

$ws.Range("T15").Value = 6.3817
$ws.Range("U15").Value = 6381471.8322

$ws.Range("T16").Value = 6.3817
$ws.Range("U16").Value = 12762943.3922

$ws.Range("T17").Value = 6.3817
$ws.Range("U17").Value = 15953678.8999

$ws.Range("T18").Value = 6.3817
$ws.Range("U18").Value = 1749448.5212

$ws.Range("T19").Value = 6.3817
$ws.Range("U19").Value = 3190735.78

$ws.Range("T31").Value = 2.1819
$ws.Range("U31").Value = 3491.0509
$ws.Range("V31").Value = 3.4818

$ws.Range("T34").Value = 84.8843
$ws.Range("U34").Value = 848818.9044999999
$ws.Range("V34").Value = 0.8065

$ws.Range("T35").Value = 82.98909999999999
$ws.Range("U35").Value = 331956.3307
$ws.Range("V35").Value = 0.8065

$ws.Range("T36").Value = 82.98909999999999
$ws.Range("U36").Value = 331956.3307
$ws.Range("V36").Value = 0.8065

$ws.Range("T37").Value = 671.5915
$ws.Range("U37").Value = 134317.2859
$ws.Range("V37").Value = 8.2247

$ws.Range("T38").Value = 79873.7458
$ws.Range("U38").Value = 1916969.9
$ws.Range("V38").Value = 984.2729

$ws.Range("T43").Value = 262171.2558
$ws.Range("U43").Value = 1310856.2557
$ws.Range("V43").Value = 2655.7896
